# Fixing typos in documentation
# - B23 ("?") on Bob's sheet becomes the numeric hours value 2.75
# - C23's task description text gets "updating README" inserted
# - B27's SUM formula recalculates automatically once B23 is numeric

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bob")

# Replace the placeholder "?" hours entry with the actual numeric hours.
$ws.Range("B23").Value = 2.75

# Fix the typo / add missing detail in the task description.
$ws.Range("C23").Value = "Finishing documentation, video, updating README, arranging demo with reviewers"

$wb.Application.CalculateFullRebuild()
